$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column (D) keeps its original text storage type.
# Several prices are numeric-looking strings (e.g. "241.55") which Excel
# would otherwise auto-convert to a floating point Number on assignment,
# so the column is pre-formatted as Text before the new values are written.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Cells.Item(2, 4).Value = '29.431.85'
$ws.Cells.Item(2, 5).Value = '  -0.49%  '

$ws.Cells.Item(3, 4).Value = '1.850.14'
$ws.Cells.Item(3, 5).Value = '  -0.23%  '

$ws.Cells.Item(4, 4).Value = '0.9992'
$ws.Cells.Item(4, 5).Value = '  -0.05%  '

$ws.Cells.Item(5, 4).Value = '241.55'
$ws.Cells.Item(5, 5).Value = '  -0.88%  '

$ws.Cells.Item(6, 4).Value = '0.6335'
$ws.Cells.Item(6, 5).Value = '  -1.00%  '

$ws.Cells.Item(7, 5).Value = '  +0.00%  '

$ws.Cells.Item(8, 4).Value = '6.405.39'
$ws.Cells.Item(8, 5).Value = '  +242.48%  '

$ws.Cells.Item(9, 4).Value = '6.431.62'
$ws.Cells.Item(9, 5).Value = '  +198.42%  '

$ws.Cells.Item(10, 4).Value = '0.07589'
$ws.Cells.Item(10, 5).Value = '  +0.45%  '

$ws.Cells.Item(11, 4).Value = '0.2967'
$ws.Cells.Item(11, 5).Value = '  -1.28%  '

$ws.Cells.Item(12, 4).Value = '24.76'
$ws.Cells.Item(12, 5).Value = '  +1.21%  '

$ws.Cells.Item(13, 4).Value = '0.07731'
$ws.Cells.Item(13, 5).Value = '  +0.82%  '

$ws.Cells.Item(14, 4).Value = '5.010'
$ws.Cells.Item(14, 5).Value = '  -0.86%  '

$ws.Cells.Item(15, 4).Value = '0.6842'
$ws.Cells.Item(15, 5).Value = '  -0.85%  '

$ws.Cells.Item(16, 5).Value = '  -1.16%  '

$ws.Cells.Item(17, 4).Value = '0.000009937'
$ws.Cells.Item(17, 5).Value = '  +3.13%  '

$ws.Cells.Item(18, 4).Value = '6.161'
$ws.Cells.Item(18, 5).Value = '  -1.81%  '

$ws.Cells.Item(19, 4).Value = '29.467.67'
$ws.Cells.Item(19, 5).Value = '  -0.59%  '

$ws.Cells.Item(20, 4).Value = '231.96'
$ws.Cells.Item(20, 5).Value = '  -2.90%  '

$ws.Cells.Item(21, 4).Value = '12.45'
$ws.Cells.Item(21, 5).Value = '  -1.42%  '

$ws.Cells.Item(22, 4).Value = '1.000'
$ws.Cells.Item(22, 5).Value = '  +0.05%  '

$ws.Cells.Item(23, 4).Value = '7.578'

$ws.Cells.Item(24, 4).Value = '1.001'
$ws.Cells.Item(24, 5).Value = '  +0.01%  '

$ws.Cells.Item(25, 4).Value = '6.590.75'
$ws.Cells.Item(25, 5).Value = '  +219.72%  '

$ws.Cells.Item(26, 2).Value = 'FraxShare'
$ws.Cells.Item(26, 3).Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Cells.Item(26, 4).Value = '16.18'
$ws.Cells.Item(26, 5).Value = '  +163.97%  '

$ws.Cells.Item(27, 2).Value = 'Monero'
$ws.Cells.Item(27, 3).Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Cells.Item(27, 4).Value = '156.20'
$ws.Cells.Item(27, 5).Value = '  -0.67%  '

$ws.Cells.Item(28, 2).Value = 'Stellar'
$ws.Cells.Item(28, 3).Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Cells.Item(28, 4).Value = '0.1396'
$ws.Cells.Item(28, 5).Value = '  -0.60%  '

$ws.Cells.Item(29, 2).Value = 'Cosmos'
$ws.Cells.Item(29, 3).Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Cells.Item(29, 4).Value = '8.428'
$ws.Cells.Item(29, 5).Value = '  -0.96%  '

$ws.Cells.Item(30, 4).Value = '17.72'
$ws.Cells.Item(30, 5).Value = '  -0.69%  '

$ws.Cells.Item(31, 2).Value = 'Frax'
$ws.Cells.Item(31, 3).Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Cells.Item(31, 4).Value = '2.611'
$ws.Cells.Item(31, 5).Value = '  +161.99%  '

$ws.Cells.Item(32, 2).Value = 'PancakeSwap'
$ws.Cells.Item(32, 3).Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Cells.Item(32, 4).Value = '1.476'
$ws.Cells.Item(32, 5).Value = '  -1.04%  '

$ws.Cells.Item(33, 4).Value = '0.05709'
$ws.Cells.Item(33, 5).Value = '  -3.36%  '

$ws.Cells.Item(34, 4).Value = '1.261'
$ws.Cells.Item(34, 5).Value = '  -1.80%  '

$ws.Cells.Item(35, 4).Value = '4.139'
$ws.Cells.Item(35, 5).Value = '  -0.22%  '

$ws.Cells.Item(36, 4).Value = '4.037'
$ws.Cells.Item(36, 5).Value = '  -1.21%  '

$ws.Cells.Item(37, 5).Value = '  -4.07%  '

$ws.Cells.Item(38, 4).Value = '1.158'
$ws.Cells.Item(38, 5).Value = '  -2.09%  '

$ws.Cells.Item(39, 4).Value = '0.7200'
$ws.Cells.Item(39, 5).Value = '  -1.01%  '

$ws.Cells.Item(40, 4).Value = '2.601'
$ws.Cells.Item(40, 5).Value = '  +0.00%  '

$ws.Cells.Item(41, 4).Value = '1.251.98'
$ws.Cells.Item(41, 5).Value = '  +3.19%  '

$ws.Cells.Item(42, 4).Value = '2.813'
$ws.Cells.Item(42, 5).Value = '  +0.40%  '

$ws.Cells.Item(43, 4).Value = '0.01811'
$ws.Cells.Item(43, 5).Value = '  +1.85%  '

$ws.Cells.Item(44, 2).Value = 'Illuvium'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/WvoRtQhzN+illuvium-ilv'
$ws.Cells.Item(44, 4).Value = '169.74'
$ws.Cells.Item(44, 5).Value = '  +258.16%  '

$ws.Cells.Item(45, 2).Value = 'TrustWalletToken'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Cells.Item(45, 4).Value = '0.9034'
$ws.Cells.Item(45, 5).Value = '  -1.43%  '

$ws.Cells.Item(46, 2).Value = 'BoneShibaSwap'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/V2gPy4UsR+boneshibaswap-bone'
$ws.Cells.Item(46, 4).Value = '5.574'
$ws.Cells.Item(46, 5).Value = '  +244.79%  '

$ws.Cells.Item(47, 2).Value = 'PaxDollar'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Cells.Item(47, 4).Value = '1.000'
$ws.Cells.Item(47, 5).Value = '  +0.03%  '

$ws.Cells.Item(48, 2).Value = 'Quant'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Cells.Item(48, 4).Value = '101.91'
$ws.Cells.Item(48, 5).Value = '  -0.17%  '

$ws.Cells.Item(49, 2).Value = 'Aave'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Cells.Item(49, 4).Value = '66.49'
$ws.Cells.Item(49, 5).Value = '  -1.41%  '

$ws.Cells.Item(50, 4).Value = '7.089'
$ws.Cells.Item(50, 5).Value = '  -5.45%  '

$ws.Cells.Item(51, 4).Value = '9.183'
$ws.Cells.Item(51, 5).Value = '  +0.28%  '
